$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Save" header in H1, re-using the same style as the existing
# header cells (bold, centered, thin border) by copying G1's formatting.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Save"

# Save flag values for rows 2-43 (1 = saved, 0 = not saved)
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 0
    23 = 1
    24 = 0
    25 = 0
    26 = 0
    27 = 1
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 1
    35 = 0
    36 = 0
    37 = 0
    38 = 1
    39 = 0
    40 = 0
    41 = 0
    42 = 1
    43 = 1
}

foreach ($r in $saveValues.Keys) {
    $ws.Cells.Item($r, 8).Value = $saveValues[$r]
}
